# Update crypto price/volume figures per latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    # Preserve the cells existing style while forcing the new value to
    # be stored as literal text (matches the workbooks inlineStr cells -
    # without this, Excel auto-coerces numeric-looking strings to numbers).
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '25.913.92'
Set-TextValue $ws.Range('E2') '  -1.38%  '
Set-TextValue $ws.Range('D3') '1.636.78'
Set-TextValue $ws.Range('E3') '  -0.65%  '
Set-TextValue $ws.Range('E4') '  +0.11%  '
Set-TextValue $ws.Range('D5') '215.49'
Set-TextValue $ws.Range('E5') '  -0.73%  '
Set-TextValue $ws.Range('E6') '  +0.16%  '
Set-TextValue $ws.Range('E7') '  +0.10%  '
Set-TextValue $ws.Range('E8') '  -0.92%  '
Set-TextValue $ws.Range('E10') '  -2.00%  '
Set-TextValue $ws.Range('D11') '0.0793'
Set-TextValue $ws.Range('E11') '  -0.05%  '
Set-TextValue $ws.Range('D12') '1.863.27'
Set-TextValue $ws.Range('E12') '  -0.65%  '
Set-TextValue $ws.Range('E13') '  -0.49%  '
Set-TextValue $ws.Range('D14') '1.638.40'
Set-TextValue $ws.Range('E14') '  +0.31%  '
Set-TextValue $ws.Range('D15') '0.544'
Set-TextValue $ws.Range('E15') '  -0.62%  '
Set-TextValue $ws.Range('E16') '  -0.54%  '
Set-TextValue $ws.Range('D17') '62.85'
Set-TextValue $ws.Range('E17') '  -0.88%  '
Set-TextValue $ws.Range('D18') '25.927.78'
Set-TextValue $ws.Range('E18') '  -1.31%  '
Set-TextValue $ws.Range('E19') '  +0.17%  '
Set-TextValue $ws.Range('D20') '192.76'
Set-TextValue $ws.Range('E20') '  -1.49%  '
Set-TextValue $ws.Range('E21') '  -2.01%  '
Set-TextValue $ws.Range('E22') '  -1.70%  '
Set-TextValue $ws.Range('E23') '  -0.82%  '
Set-TextValue $ws.Range('E24') '  +4.60%  '
Set-TextValue $ws.Range('E25') '  +0.49%  '
Set-TextValue $ws.Range('D26') '143.29'
Set-TextValue $ws.Range('E26') '  -0.23%  '
Set-TextValue $ws.Range('E27') '  +0.05%  '
Set-TextValue $ws.Range('D28') '6.88'
Set-TextValue $ws.Range('E28') '  -1.02%  '
Set-TextValue $ws.Range('E29') '  -0.63%  '
Set-TextValue $ws.Range('D30') '1.25'
Set-TextValue $ws.Range('E30') '  -0.47%  '
Set-TextValue $ws.Range('E31') '  -0.72%  '
Set-TextValue $ws.Range('E32') '  -2.37%  '
Set-TextValue $ws.Range('E33') '  -0.36%  '
Set-TextValue $ws.Range('E34') '  -4.40%  '
Set-TextValue $ws.Range('E35') '  +1.50%  '
Set-TextValue $ws.Range('E36') '  -1.51%  '
Set-TextValue $ws.Range('D37') '1.133.80'
Set-TextValue $ws.Range('E37') '  -0.34%  '
Set-TextValue $ws.Range('E38') '  -1.86%  '
Set-TextValue $ws.Range('E39') '  -1.67%  '
Set-TextValue $ws.Range('E40') '  -0.51%  '
Set-TextValue $ws.Range('D41') '5.49'
Set-TextValue $ws.Range('E41') '  -0.75%  '
Set-TextValue $ws.Range('D42') '99.42'
Set-TextValue $ws.Range('E42') '  -1.18%  '
Set-TextValue $ws.Range('E43') '  -0.52%  '
Set-TextValue $ws.Range('D44') '1.773.17'
Set-TextValue $ws.Range('D45') '0.0₆0114'
Set-TextValue $ws.Range('E45') '  +2.07%  '
Set-TextValue $ws.Range('D46') '56.51'
Set-TextValue $ws.Range('E46') '  -1.22%  '
Set-TextValue $ws.Range('E47') '  +2.19%  '
Set-TextValue $ws.Range('E48') '  +0.08%  '
Set-TextValue $ws.Range('D49') '7.67'
Set-TextValue $ws.Range('E49') '  +0.04%  '
Set-TextValue $ws.Range('E50') '  -0.94%  '
Set-TextValue $ws.Range('D51') '0.0958'
Set-TextValue $ws.Range('E51') '  -1.42%  '
